# Apply the 7 May 2023 cryptos-list refresh:
#  - rows 2-23 get updated Price (D) / Volume(1h) (E) figures
#  - a new coin (WrappedliquidstakedEther2.0) is inserted at row 24,
#    pushing every coin from the old row 24 (Toncoin) through the old
#    row 50 (Quant) down by one row, each keeping its own refreshed
#    Price/Volume figures; the old last row (PaxDollar) drops off the
#    bottom of the A1:E51 range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Table of per-row column updates. Each entry is Row -> @{ Col = NewValue }.
# Using Range(...).NumberFormat = "@" before writing keeps values such as
# "28.999.50" or "0.3823" stored as literal text instead of being
# reinterpreted/rounded as numbers by Excel.
$rowUpdates = [ordered]@{
    2 = @{ 'D'='28.999.50'; 'E'='  -1.69%  ' }
    3 = @{ 'D'='1.909.11'; 'E'='  -3.17%  ' }
    4 = @{ 'E'='  +0.06%  ' }
    5 = @{ 'D'='324.19'; 'E'='  -0.76%  ' }
    6 = @{ 'E'='  -0.22%  ' }
    7 = @{ 'E'='  -1.31%  ' }
    8 = @{ 'D'='0.3823'; 'E'='  -2.08%  ' }
    9 = @{ 'E'='  -2.94%  ' }
    10 = @{ 'D'='0.9796'; 'E'='  -0.79%  ' }
    11 = @{ 'D'='22.07'; 'E'='  -3.10%  ' }
    12 = @{ 'D'='1.889.73'; 'E'='  -4.29%  ' }
    13 = @{ 'D'='5.675'; 'E'='  -2.16%  ' }
    14 = @{ 'D'='6.933'; 'E'='  -3.23%  ' }
    15 = @{ 'D'='0.07021'; 'E'='  -1.09%  ' }
    16 = @{ 'E'='  -0.19%  ' }
    17 = @{ 'D'='83.88'; 'E'='  -4.36%  ' }
    18 = @{ 'D'='0.000009454'; 'E'='  -4.74%  ' }
    19 = @{ 'E'='  -3.05%  ' }
    21 = @{ 'D'='28.952.53'; 'E'='  -1.97%  ' }
    22 = @{ 'E'='  -3.88%  ' }
    23 = @{ 'E'='  -2.24%  ' }
    24 = @{ 'B'='WrappedliquidstakedEther2.0'; 'C'='https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'; 'D'='2.115.89'; 'E'='  -4.61%  ' }
    25 = @{ 'B'='Toncoin'; 'C'='https://coinranking.com/coin/67YlI0K1b+toncoin-ton'; 'D'='2.093'; 'E'='  -0.47%  ' }
    26 = @{ 'B'='Monero'; 'C'='https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; 'D'='158.40'; 'E'='  -0.12%  ' }
    27 = @{ 'B'='EthereumClassic'; 'C'='https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; 'D'='19.07'; 'E'='  -2.09%  ' }
    28 = @{ 'B'='InternetComputer(DFINITY)'; 'C'='https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; 'D'='5.666'; 'E'='  -2.00%  ' }
    29 = @{ 'B'='BitcoinCash'; 'C'='https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'; 'D'='117.52'; 'E'='  -1.76%  ' }
    30 = @{ 'B'='LidoDAOToken'; 'C'='https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'; 'D'='1.850'; 'E'='  -1.76%  ' }
    31 = @{ 'B'='Stellar'; 'C'='https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; 'D'='0.09268'; 'E'='  -1.57%  ' }
    32 = @{ 'B'='ImmutableX'; 'C'='https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; 'D'='0.8644'; 'E'='  -1.26%  ' }
    33 = @{ 'B'='Filecoin'; 'C'='https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; 'D'='5.069'; 'E'='  -2.95%  ' }
    34 = @{ 'B'='ARBITRUM'; 'C'='https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'; 'D'='1.248'; 'E'='  -5.41%  ' }
    35 = @{ 'B'='HuobiToken'; 'C'='https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'; 'D'='3.025'; 'E'='  -3.30%  ' }
    36 = @{ 'B'='Hedera'; 'C'='https://coinranking.com/coin/jad286TjB+hedera-hbar'; 'D'='0.05742'; 'E'='  -1.02%  ' }
    37 = @{ 'B'='TrustWalletToken'; 'C'='https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; 'D'='1.156'; 'E'='  -0.62%  ' }
    38 = @{ 'B'='Frax'; 'C'='https://coinranking.com/coin/KfWtaeV1W+frax-frax'; 'D'='1.001'; 'E'='  -0.24%  ' }
    39 = @{ 'B'='VeChain'; 'C'='https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; 'D'='0.02042'; 'E'='  -2.86%  ' }
    40 = @{ 'B'='TheSandbox'; 'C'='https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'; 'D'='0.5506'; 'E'='  -3.43%  ' }
    41 = @{ 'B'='FraxShare'; 'C'='https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'; 'D'='7.407'; 'E'='  -3.79%  ' }
    42 = @{ 'B'='Algorand'; 'C'='https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'; 'D'='0.1757'; 'E'='  -2.04%  ' }
    43 = @{ 'B'='MXToken'; 'C'='https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'; 'D'='2.851'; 'E'='  +3.54%  ' }
    44 = @{ 'B'='Aptos'; 'C'='https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; 'D'='9.317'; 'E'='  -3.12%  ' }
    45 = @{ 'B'='Decentraland'; 'C'='https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'; 'D'='0.5177'; 'E'='  -2.75%  ' }
    46 = @{ 'B'='EnergySwap'; 'C'='https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; 'D'='11.31'; 'E'='  -2.99%  ' }
    47 = @{ 'B'='Cronos'; 'C'='https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'; 'D'='0.06836'; 'E'='  -1.11%  ' }
    48 = @{ 'B'='PEPE'; 'C'='https://coinranking.com/coin/03WI8NQPF+pepe-pepe'; 'D'='0.000002600'; 'E'='  -6.95%  ' }
    49 = @{ 'B'='RenderToken'; 'C'='https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; 'D'='2.048'; 'E'='  -4.69%  ' }
    50 = @{ 'B'='Quant'; 'C'='https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'; 'D'='110.95'; 'E'='  -1.94%  ' }
    51 = @{ 'B'='NEARProtocol'; 'C'='https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'; 'D'='1.780'; 'E'='  -2.61%  ' }
}

foreach ($rowNum in $rowUpdates.Keys) {
    $cols = $rowUpdates[$rowNum]
    foreach ($col in $cols.Keys) {
        $cell = $ws.Range("$col$rowNum")
        $cell.NumberFormat = "@"
        $cell.Value = $cols[$col]
        # Restore the plain "Normal" style so the cell is left exactly as it
        # started (no explicit style index), matching every other untouched
        # data cell in the sheet; the text-literal value written above is
        # unaffected by this.
        $cell.Style = "Normal"
    }
}

